$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: change C3 from "RPL" to "tet", and D3 from 3 to 1
$ws.Range("A3").Value = "gambar_2.jpg"
$ws.Range("B3").Value = "bebas"
$ws.Range("C3").Value = "tet"
$ws.Range("D3").Value = 1

# Rows 4-10: new rows with gambar_2.jpg / bebas / TET / 1
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "gambar_2.jpg"
    $ws.Cells.Item($r, 2).Value = "bebas"
    $ws.Cells.Item($r, 3).Value = "TET"
    $ws.Cells.Item($r, 4).Value = 1
}

# Update selection to D3
$ws.Range("D3").Select()
